$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 4. B4 ("Индикатор" row): update indicator name to the longer, more specific description
$ws.Range("B4").Value = "4.3.1.1. Доля молодежи от 15 до 24 лет, обучающихся в системе начального профессионального образования, среднего профессионального и высшего профессионального образования к численности населения соответствующего возраста"

# 5. B10 ("Сайт организации" row): update site url and add hyperlink
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Hyperlinks.Add($ws.Range("B10"), "http://www.stat.gov.kg")

# 6. Update selection to B10
$ws.Range("B10").Select()
